# Insert a new weekly price record for "Zanahoria" (Agrícola del Norte S.A.
# de Arica) as row 226 of Sheet1, pushing the existing rows 226-241 down to
# 227-242 (so the sheet's used range grows from A1:R241 to A1:R242).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 226..241 down to 227..242, leaving a blank row 226 in place.
$ws.Rows.Item(226).Insert()

# Populate the new row 226 with the new weekly record.
$ws.Cells.Item(226, 1).Value  = 1
$ws.Cells.Item(226, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(226, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(226, 4).Value  = 44610
$ws.Cells.Item(226, 5).Value  = 15
$ws.Cells.Item(226, 6).Value  = 100114013
$ws.Cells.Item(226, 7).Value  = "Zanahoria"
$ws.Cells.Item(226, 8).Value  = "Sin especificar"
$ws.Cells.Item(226, 9).Value  = "Primera"
$ws.Cells.Item(226, 10).Value = 70
$ws.Cells.Item(226, 11).Value = 20000
$ws.Cells.Item(226, 12).Value = 21000
$ws.Cells.Item(226, 13).Value = 20500
$ws.Cells.Item(226, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(226, 15).Value = "Valle de Camiña"
$ws.Cells.Item(226, 16).Value = 820
$ws.Cells.Item(226, 17).Value = 25
$ws.Cells.Item(226, 18).Value = "Hortaliza"
